$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("testcases")
$wsFail = $wb.Worksheets.Item("failing testcases")

# ---------------------------------------------------------------------------
# 1. Add the new "failing testcases" rows / column (this also grows the
#    shared-string table with the 14 new strings, in original author order).
# ---------------------------------------------------------------------------

# New column header
$wsFail.Range("D1").Value = "Status"

# Existing rows 7-11 (new testcases)
$wsFail.Range("A7").Value = "characters played by Cher"
$wsFail.Range("B7").Value = "error message"
$wsFail.Range("C7").Value = "casing"

$wsFail.Range("A8").Value = "characters and movies played by Cher"
$wsFail.Range("B8").Value = "list too long"
$wsFail.Range("C8").Value = "list for characters played by Cher is correct"

$wsFail.Range("A9").Value = "poster for Gone with the Wind"
$wsFail.Range("B9").Value = "shows poster for Gone Fishin"
$wsFail.Range("C9").Value = "too loose matching"

$wsFail.Range("A10").Value = "overall casing"
$wsFail.Range("B10").Value = "inconsistent results"

$wsFail.Range("A11").Value = "director of GoldenEye"
$wsFail.Range("B11").Value = "error message"
$wsFail.Range("C11").Value = "same for any query with condition_col in movies and ranked_col in a child table"

# New column D width
$wsFail.Columns.Item(4).ColumnWidth = 33.7

# ---------------------------------------------------------------------------
# 2. Row-height touch-ups on "testcases" sheet.
# ---------------------------------------------------------------------------
$wsTest.Rows.Item(1).RowHeight = 46.5
$wsTest.Rows.Item(16).RowHeight = 31
$wsTest.Rows.Item(17).RowHeight = 31
$wsTest.Rows.Item(21).RowHeight = 31
$wsTest.Rows.Item(22).RowHeight = 124
$wsTest.Rows.Item(23).RowHeight = 108.5
$wsTest.Rows.Item(24).RowHeight = 108.5

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping.
# ---------------------------------------------------------------------------
$wsTest.Range("B7").Select() | Out-Null

$wsFail.Activate()
$wsFail.Range("B17").Select() | Out-Null
